$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27
$ws.Range("B4").Value = 11596.53993302689
$ws.Range("B5").Value = 15785.03079000041
$ws.Range("B6").Value = 1170.388941922528
$ws.Range("B7").Value = 1056.008723141073
$ws.Range("B8").Value = 20217.59999999949
$ws.Range("B9").Value = 3738.737584105171
$ws.Range("B10").Value = 193226.2239903583
$ws.Range("B11").Value = 0.1176160567795326
$ws.Range("B12").Value = 0.3541473878820095
$ws.Range("B13").Value = 0.3500000000000018
$ws.Range("B14").Value = 0.9927538805401921
$ws.Range("B15").Value = 0.5510543047759427
